$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp string (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 9 de Abril de 2020 a las 09:22"

# --- Turquia (row 19): totals / new-cases / active updated ---
$ws.Range("B19").Value = 12981
$ws.Range("C19").Value = 39
$ws.Range("E19").Value = 8196

# --- Polonia (row 31): active / recovered updated ---
$ws.Range("D31").Value = 284
$ws.Range("E31").Value = 4762

# --- Moldavia (row 64): active / recovered updated ---
$ws.Range("D64").Value = 43
$ws.Range("E64").Value = 1103

# --- Armenia (row 70): totals / new-cases / active / recovered / critical / deaths-today updated ---
$ws.Range("B70").Value = 921
$ws.Range("C70").Value = 40
$ws.Range("D70").Value = 138
$ws.Range("E70").Value = 773
$ws.Range("G70").Value = 1
$ws.Range("H70").Value = 10

# --- Bosnia y Herzegovina (row 72): totals / new-cases / recovered updated ---
$ws.Range("B72").Value = 839
$ws.Range("C72").Value = 35
$ws.Range("E72").Value = 709

# --- Letonia (row 82): totals / new-cases / recovered / critical updated ---
$ws.Range("B82").Value = 589
$ws.Range("C82").Value = 12
$ws.Range("E82").Value = 571
$ws.Range("F82").Value = 3

# --- Montenegro / Vietnam swap order + updated figures ---
# Row 108 used to be Vietnam, row 109 used to be Montenegro; the countries
# now swap places (Montenegro moves above Vietnam) with refreshed data.
$ws.Range("A108").Value = "Montenegro"
$ws.Range("B108").Value = 252
$ws.Range("C108").Value = 4
$ws.Range("D108").Value = 4
$ws.Range("E108").Value = 246
$ws.Range("F108").Value = 7
$ws.Range("G108").Value = 0
$ws.Range("H108").Value = 2

$ws.Range("A109").Value = "Vietnam"
$ws.Range("B109").Value = 251
$ws.Range("C109").Value = 0
$ws.Range("D109").Value = 128
$ws.Range("E109").Value = 123
$ws.Range("F109").Value = 8
$ws.Range("G109").Value = 0
$ws.Range("H109").Value = 0

# --- Georgia (row 112): totals / new-cases / recovered updated ---
$ws.Range("B112").Value = 214
$ws.Range("C112").Value = 3
$ws.Range("E112").Value = 161

# --- Islas Feroe (row 115): active / recovered updated ---
$ws.Range("D115").Value = 136
$ws.Range("E115").Value = 48

# --- Surinam (row 188): active / recovered updated ---
$ws.Range("D188").Value = 4
$ws.Range("E188").Value = 5
